# Weekly update: insert 3 new rows of Chirimoya price data (date 45202)
# at the top of the data block (rows 202-204), pushing the existing
# rows 202-305 down to 205-308.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 202:305 down by inserting 3 new blank rows at 202.
$ws.Rows("202:204").Insert()

# The columns that stay constant across these three new rows (copied
# from the surrounding "Chirimoya" block pattern).
$commonA = 8
$commonB = "Terminal La Palmera de La Serena"
$commonC = "Coquimbo"
$commonE = 4
$commonF = "Fruta"
$commonG = 100107
$commonH = "Otros"
$commonI = 100107002
$commonJ = "Chirimoya"
$commonK = "Cultivar IV Región"

function Set-ChirimoyaRow {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Unidad,
        [string]$Origen,
        [double]$PrecioKg,
        [double]$KgUnidad
    )

    $ws.Cells.Item($Row, 1).Value = $commonA
    $ws.Cells.Item($Row, 2).Value = $commonB
    $ws.Cells.Item($Row, 3).Value = $commonC
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = $commonE
    $ws.Cells.Item($Row, 6).Value = $commonF
    $ws.Cells.Item($Row, 7).Value = $commonG
    $ws.Cells.Item($Row, 8).Value = $commonH
    $ws.Cells.Item($Row, 9).Value = $commonI
    $ws.Cells.Item($Row, 10).Value = $commonJ
    $ws.Cells.Item($Row, 11).Value = $commonK
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

# Row 202: Especial, Provincia de Limarí, $/bandeja 10 kilos
Set-ChirimoyaRow -Row 202 -Fecha 45202 -Calidad "Especial" -Volumen 240 `
    -PrecioMin 22500 -PrecioMax 23000 -PrecioProm 22750 `
    -Unidad "$/bandeja 10 kilos" -Origen "Provincia de Limarí" `
    -PrecioKg 2275 -KgUnidad 10

# Row 203: Primera, Provincia de Limarí, $/bandeja 10 kilos
Set-ChirimoyaRow -Row 203 -Fecha 45202 -Calidad "Primera" -Volumen 200 `
    -PrecioMin 20500 -PrecioMax 21000 -PrecioProm 20750 `
    -Unidad "$/bandeja 10 kilos" -Origen "Provincia de Limarí" `
    -PrecioKg 2075 -KgUnidad 10

# Row 204: Segunda, Provincia de Limarí, $/bandeja 10 kilos
Set-ChirimoyaRow -Row 204 -Fecha 45202 -Calidad "Segunda" -Volumen 160 `
    -PrecioMin 17500 -PrecioMax 18000 -PrecioProm 17750 `
    -Unidad "$/bandeja 10 kilos" -Origen "Provincia de Limarí" `
    -PrecioKg 1775 -KgUnidad 10

# Ensure date cells use the expected date/time style (copied automatically
# on row insert, but set explicitly here for safety).
$ws.Range("D202:D204").NumberFormat = "YYYY-MM-DD HH:MM:SS"
